$d = $word.ActiveDocument

# Locate the paragraph holding the old, centered/italic "Notice u/s 94 BNSS"
# heading and remove the whole paragraph (including its paragraph mark).
$noticePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Notice u/s 94 BNSS" + [char]13) {
        $noticePara = $p
        break
    }
}
$noticeIndex = $noticePara.Index
$noticePara.Range.Delete()

# After the delete, the paragraph that now sits at the old index is the blank
# "Body A" paragraph that used to directly follow the removed heading (the
# first of the two blank paragraphs above "To,"). Insert a brand new, empty
# paragraph right after it so it is cleanly isolated between two real
# paragraph marks - this leaves both surrounding blank paragraphs completely
# untouched, instead of merging new content into one of them.
$blankAbove = $d.Paragraphs.Item($noticeIndex)
$blankAbove.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($noticeIndex + 1)

# Overwrite that isolated (still empty) paragraph with the exact target
# markup: no paragraph style, justified, single bold+underlined run reading
# "Notice u/s 94 BNSS, 2023".
$xmlPackage = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Notice u/s 94 BNSS, 2023</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$null = $newPara.Range.InsertXML($xmlPackage)
